# Generate Report for Handback
# Update the "generate date" / handoff-handback timestamp cells to reflect
# the newly generated report times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G2).
$wsOverview.Range("G2").Value = "2017-02-21 04:54:04"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (L2)
$wsZhCn.Range("H2").Value = "2017-02-21 04:53:45"
$wsZhCn.Range("L2").Value = "2017-02-21 04:54:39"

# de-de sheet: "Correspond Handoff Datetime" (H2) shares the same text value
# as the Overview sheet's "Latest HO Xliff Generate Date" (G2), and
# "Correspond Handback DateTime" (L2).
$wsDeDe.Range("H2").Value = "2017-02-21 04:54:04"
$wsDeDe.Range("L2").Value = "2017-02-21 04:55:02"
